$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New close values appended for rows 878-920 (A: index 876-918, B: close value)
$bValues = @(0.88049, 0.8792, 0.88195, 0.8762799999999999, 0.87093, 0.84612, 0.85017, 0.85017, 0.85017, 0.85439, 0.84612, 0.84206, 0.82536, 0.83363, 0.83785, 0.8226, 0.83104, 0.82406, 0.81158, 0.79796, 0.76342, 0.74964, 0.79796, 0.74964, 0.75661, 0.76893, 0.74282, 0.70829, 0.70829, 0.66694, 0.68623, 0.6945, 0.70829, 0.70569, 0.72077, 0.72904, 0.76066, 0.70699, 0.6684, 0.69726, 0.70569, 0.73326, 0.72207)

$startRow = 878
$startIndex = 876

# Copy the style (bold/border/alignment) from the last existing data row in column A
# so the new index cells match the existing formatting (style id 1).
$styleSource = $ws.Range("A877")

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $startRow + $i
    $styleSource.Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value2 = $startIndex + $i
    $ws.Cells.Item($row, 2).Value2 = $bValues[$i]
}

$excel.CutCopyMode = 0
